$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.086.12"
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").Value = "2.369.78"
$ws.Range("E3").Value = "  +1.27%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'303.76"
$ws.Range("E5").Value = "  +0.27%  "
$ws.Range("D6").Value = "'95.93"
$ws.Range("E6").Value = "  +0.98%  "
$ws.Range("D7").Value = "'0.504"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "'0.484"
$ws.Range("E9").Value = "  -2.49%  "
$ws.Range("D10").Value = "'34.43"
$ws.Range("E10").Value = "  +0.90%  "
$ws.Range("E11").Value = "  +4.24%  "
$ws.Range("E12").Value = "  +0.63%  "
$ws.Range("D13").Value = "'18.32"
$ws.Range("E13").Value = "  -2.10%  "
$ws.Range("D14").Value = "'6.81"
$ws.Range("E14").Value = "  +0.73%  "
$ws.Range("D15").Value = "2.734.93"
$ws.Range("E15").Value = "  +1.11%  "
$ws.Range("D16").Value = "2.359.66"
$ws.Range("E16").Value = "  -0.23%  "
$ws.Range("D17").Value = "'0.804"
$ws.Range("E17").Value = "  +0.55%  "
$ws.Range("D18").Value = "43.098.36"
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("D19").Value = "'11.97"
$ws.Range("E19").Value = "  -1.78%  "
$ws.Range("E20").Value = "  +1.17%  "
$ws.Range("D22").Value = "'67.96"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").Value = "'235.46"
$ws.Range("E23").Value = "  -0.25%  "
$ws.Range("E24").Value = "  -0.41%  "
$ws.Range("D25").Value = "'2.45"
$ws.Range("E25").Value = "  +1.48%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").Value = "'24.50"
$ws.Range("E27").Value = "  -0.55%  "
$ws.Range("E28").Value = "  +0.59%  "
$ws.Range("D29").Value = "'9.37"
$ws.Range("E29").Value = "  +2.19%  "
$ws.Range("D30").Value = "'32.03"
$ws.Range("E30").Value = "  +1.48%  "
$ws.Range("E31").Value = "  -0.17%  "
$ws.Range("D32").Value = "'5.05"
$ws.Range("E32").Value = "  +0.73%  "
$ws.Range("E33").Value = "  +10.75%  "
$ws.Range("D34").Value = "'17.82"
$ws.Range("E34").Value = "  +3.10%  "
$ws.Range("D35").Value = "'0.0737"
$ws.Range("E35").Value = "  +1.45%  "
$ws.Range("D36").Value = "'128.80"
$ws.Range("E36").Value = "  +12.91%  "
$ws.Range("D37").Value = "'1.83"
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("E38").Value = "  +3.41%  "
$ws.Range("D39").Value = "'4.33"
$ws.Range("E39").Value = "  -1.48%  "
$ws.Range("E40").Value = "  -3.20%  "
$ws.Range("E41").Value = "  -0.68%  "
$ws.Range("D42").Value = "'21.19"
$ws.Range("E42").Value = "  -4.95%  "
$ws.Range("D43").Value = "1.929.42"
$ws.Range("E43").Value = "  -0.35%  "
$ws.Range("D44").Value = "'0.0279"
$ws.Range("E44").Value = "  -1.10%  "
$ws.Range("E45").Value = "  +1.64%  "
$ws.Range("D46").Value = "'2.77"
$ws.Range("E46").Value = "  +1.50%  "
$ws.Range("D47").Value = "'9.20"
$ws.Range("E47").Value = "  -8.53%  "
$ws.Range("D48").Value = "2.595.09"
$ws.Range("E48").Value = "  +1.12%  "
$ws.Range("E49").Value = "  +2.88%  "
$ws.Range("D50").Value = "'71.69"
$ws.Range("E50").Value = "  -0.64%  "
$ws.Range("D51").Value = "'51.69"
